# edit.ps1 - Applies the LOM3085.docx changes described by the commit diff.
$d = $word.ActiveDocument

function Replace-Exact([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
    }
    return $found
}

# 1) Créditos-aula: 4 -> 2
Replace-Exact "Créditos-aula: 4" "Créditos-aula: 2"

# 2) Carga horária: 60 h -> 30 h
Replace-Exact "Carga horária: 60 h" "Carga horária: 30 h"

# 3) Ativação: 01/01/2020 -> 01/01/2025
Replace-Exact "Ativação: 01/01/2020" "Ativação: 01/01/2025"

# 4) Programa resumido paragraph: drop the calorimetry/thermal-analysis clause
Replace-Exact "Técnicas de Materialografia. Calorimetria e análises térmicas de materiais." "Técnicas de Materialografia."

# 5) Programa section: remove the "2. CALORIMETRIA ..." sentence (and the manual
#    line break that precedes it) from the end of the "1. MATERIALOGRAFIA ..." run.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("2. CALORIMETRIA*FTIR).", $false, $false, $true, $false, $false, $true, 1, $false, "", 0)
if ($found2) {
    $full = $d.Range($rng2.Start - 1, $rng2.End)
    $full.Text = ""
} else {
    Write-Output "NOT FOUND: CALORIMETRIA sentence"
}

# 6) Bibliography paragraph: multiple small corrections + drop the thermal-analysis
#    reference list (trailing citations removed).
$biblioOld = "COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, São Paulo – 1974.COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.AZEVEDO, A. D.; MOTHE, C. G. Análaise Térmica de Materiais. São Paulo: ARTLIBER, 2009.BROWN, M.E. Handbook of Thermal Analysis and Calorimetry, Amsterdam: Elsevier Science, 1998.HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, New York: Wiley, 1999.HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.MULLER, A. Solidificação e Análise Térmica dos Metais. Porto Alegre: Ed. UFRGS, 2002.SPEYER, R. Thermal analysis of materials, New York: Marcel Dekker, 1994.REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982. Nondestructive Characterization of Materials. Series. Plenum Press, New York. YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994."
$biblioNew = "COLPAERT; HUBERTUS. Metalografia dos produtos siderúrgicos comuns, 3ª Edição, Editora Edgard Blücher Ltda, SãoPaulo – 1974.COUTINHO, TELMO DE AZEVEDO. Metalografia de Não-Ferrosos, Editora Edgard Blücher Ltda, São Paulo – 1980.PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.REED-HILL, R.E. Princípios de Metalurgia Física, Ed. Guanabara Dois, 1982.Nondestructive Characterization of Materials. Series. Plenum Press, New York.YACOBI, B.G.; HOLT, D.B.; KAZMERSKI, L.L. Eds. Microanalysis of Solids. Plenum Press, New York, 1994."
Replace-Exact $biblioOld $biblioNew
